# Updates the crypto price/volume table with the latest scraped values.
# (Sat Nov 18 21:52:38 UTC 2023 GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.565.24"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").Value = "1.961.09"
$ws.Range("E3").Value = "  +0.84%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("D6").Value = "'0.618"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").Value = "'58.82"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.47%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +3.16%  "

$ws.Range("E10").Value = "  -2.50%  "

$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").Value = "'22.25"
$ws.Range("D12").Style = "Normal"

$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("D14").Value = "'0.828"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.64%  "

$ws.Range("D15").Value = "'13.72"
$ws.Range("D15").Style = "Normal"

$ws.Range("E16").Value = "  +0.70%  "

$ws.Range("D17").Value = "1.958.94"
$ws.Range("E17").Value = "  +1.40%  "

$ws.Range("D18").Value = "36.498.77"
$ws.Range("E18").Value = "  +0.50%  "

$ws.Range("D19").Value = "'69.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.53%  "

$ws.Range("D20").Value = "0.0₃0859"
$ws.Range("E20").Value = "  -0.20%  "

$ws.Range("D21").Value = "'228.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").Value = "'5.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.34%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").Value = "'2.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.36%  "

$ws.Range("D25").Value = "'2.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.89%  "

# Row 26/27: Kaspa and Cosmos swap positions in the ranking, with refreshed values.
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'9.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.140"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.09%  "

$ws.Range("D28").Value = "'160.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.96%  "

$ws.Range("D29").Value = "'19.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("D30").Value = "'0.120"
$ws.Range("D30").Style = "Normal"

$ws.Range("E31").Value = "  +1.12%  "

$ws.Range("E32").Value = "  +1.40%  "

$ws.Range("D33").Value = "'0.0620"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.14%  "

$ws.Range("D34").Value = "'4.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.73%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").Value = "'2.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.94%  "

$ws.Range("D37").Value = "'3.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.02%  "

$ws.Range("E38").Value = "  -4.84%  "

$ws.Range("E39").Value = "  -0.13%  "

$ws.Range("D40").Value = "'0.0986"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("E41").Value = "  +0.97%  "

$ws.Range("E42").Value = "  +0.26%  "

$ws.Range("E43").Value = "  +0.95%  "

$ws.Range("D44").Value = "'16.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").Value = "1.365.44"
$ws.Range("E45").Value = "  +0.66%  "

$ws.Range("D46").Value = "'1.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.23%  "

$ws.Range("D47").Value = "'87.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("D48").Value = "'7.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("E49").Value = "  +0.61%  "

$ws.Range("D50").Value = "2.140.21"
$ws.Range("E50").Value = "  +0.93%  "

$ws.Range("D51").Value = "'43.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.60%  "

Write-Output "cryptos list updated"
